# Daily attendance processing - 2025-11-28 17:24:19
# Swap the order of "Recorded By" contributors in column G from
# "System, dnasr281@gmail.com" to "dnasr281@gmail.com, System"
# for every row on the "Session Analysis Results" sheet that currently
# has that exact value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

# "Recorded By" is column G (7). Used range runs to row 157.
$lastRow = 157
$updated = 0

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
        $updated = $updated + 1
    }
}

Write-Output "Updated $updated cell(s) in column G"
